$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 575.3333
$ws.Range("J33").Value = 646.75
$ws.Range("L33").Value = 646.75
$ws.Range("N33").Value = -1104.75
$ws.Range("H69").Value = 12501
$ws.Range("I69").Value = 12003
$ws.Range("J69").Value = 12750
$ws.Range("K69").Value = 36009
$ws.Range("L69").Value = 38250
$ws.Range("M69").Value = -35135
$ws.Range("N69").Value = -39998
$ws.Range("H70").Value = 2472.2222
$ws.Range("I70").Value = 1775
$ws.Range("K70").Value = 5325
$ws.Range("M70").Value = -5055
$ws.Range("H72").Value = 12501
$ws.Range("I72").Value = 12003
$ws.Range("J72").Value = 12750
$ws.Range("K72").Value = 108027
$ws.Range("L72").Value = 114750
$ws.Range("M72").Value = -103659
$ws.Range("N72").Value = -123486
$ws.Range("H73").Value = 2472.2222
$ws.Range("I73").Value = 1775
$ws.Range("K73").Value = 5325
$ws.Range("M73").Value = -4389
$ws.Range("H100").Value = 3674.625
$ws.Range("I100").Value = 3674.625
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3674.625
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -3133.625
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 7533.1665
$ws.Range("I113").Value = 7155.4443
$ws.Range("K113").Value = 7155.4443
$ws.Range("M113").Value = -3901.4443
$ws.Range("H138").Value = 2371.647
$ws.Range("I138").Value = 4164.5557
$ws.Range("K138").Value = 12493.6671
$ws.Range("M138").Value = -7353.667099999999

$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 4638.654
$ws.Range("I32").Value = 4224.24
$ws.Range("J32").Value = 14999
$ws.Range("K32").Value = 4224.24
$ws.Range("L32").Value = 14999
$ws.Range("M32").Value = -3937.24
$ws.Range("N32").Value = -15573
$ws.Range("H80").Value = 19999
$ws.Range("I80").Value = 19999
$ws.Range("K80").Value = 19999
$ws.Range("M80").Value = -19001
$ws.Range("H83").Value = 19999
$ws.Range("I83").Value = 19999
$ws.Range("K83").Value = 59997
$ws.Range("M83").Value = -55005
$ws.Range("H122").Value = 1333
$ws.Range("I122").Value = 999.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2998.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -548.5
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 1420.0769
$ws.Range("I20").Value = 1723.4445
$ws.Range("J20").Value = 737.5
$ws.Range("K20").Value = 1723.4445
$ws.Range("L20").Value = 737.5
$ws.Range("M20").Value = -1476.4445
$ws.Range("N20").Value = -1231.5
$ws.Range("H82").Value = 14432.833
$ws.Range("I82").Value = 14432.833
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 14432.833
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -14049.833
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 14432.833
$ws.Range("I85").Value = 14432.833
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 14432.833
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -13106.833
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 14837.708
$ws.Range("J86").Value = 9943.777
$ws.Range("L86").Value = 9943.777
$ws.Range("N86").Value = -12189.777
$ws.Range("H89").Value = 14837.708
$ws.Range("J89").Value = 9943.777
$ws.Range("L89").Value = 49718.885
$ws.Range("N89").Value = -60950.885
$ws.Range("H107").Value = 1096.4
$ws.Range("I107").Value = 1096.4
$ws.Range("K107").Value = 1096.4
$ws.Range("M107").Value = 823.5999999999999
$ws.Range("H134").Value = 1046.5555
$ws.Range("I134").Value = 1034.76
$ws.Range("K134").Value = 3104.28
$ws.Range("M134").Value = -569.2799999999997

$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 7732.857
$ws.Range("J31").Value = 5406.4614
$ws.Range("L31").Value = 5406.4614
$ws.Range("N31").Value = -5996.4614
$ws.Range("H34").Value = 7732.857
$ws.Range("J34").Value = 5406.4614
$ws.Range("L34").Value = 5406.4614
$ws.Range("N34").Value = -5810.4614
$ws.Range("H51").Value = 36221.285
$ws.Range("I51").Value = 33111
$ws.Range("J51").Value = 41819.8
$ws.Range("K51").Value = 33111
$ws.Range("L51").Value = 41819.8
$ws.Range("M51").Value = -32375
$ws.Range("N51").Value = -43291.8
$ws.Range("H60").Value = 22969.6
$ws.Range("I60").Value = 3198.6
$ws.Range("J60").Value = 42740.6
$ws.Range("K60").Value = 3198.6
$ws.Range("L60").Value = 42740.6
$ws.Range("M60").Value = -2687.6
$ws.Range("N60").Value = -43762.6
$ws.Range("H61").Value = 36221.285
$ws.Range("I61").Value = 33111
$ws.Range("J61").Value = 41819.8
$ws.Range("K61").Value = 33111
$ws.Range("L61").Value = 41819.8
$ws.Range("M61").Value = -32763
$ws.Range("N61").Value = -42515.8
$ws.Range("H62").Value = 4596.6
$ws.Range("I62").Value = 4596.6
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4596.6
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3972.6
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4596.6
$ws.Range("I65").Value = 4596.6
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 22983
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -19863
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 50313
$ws.Range("J74").Value = 50313
$ws.Range("L74").Value = 50313
$ws.Range("N74").Value = -52061
$ws.Range("H77").Value = 50313
$ws.Range("J77").Value = 50313
$ws.Range("L77").Value = 150939
$ws.Range("N77").Value = -159675
$ws.Range("H107").Value = 1332.8572
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H121").Value = 29073.309
$ws.Range("J121").Value = 29073.309
$ws.Range("L121").Value = 29073.309
$ws.Range("N121").Value = -31693.309
$ws.Range("H141").Value = 362981.72
$ws.Range("J141").Value = 362981.72
$ws.Range("L141").Value = 362981.72
$ws.Range("N141").Value = -373341.72

$ws = $wb.Worksheets.Item(5)
$ws.Range("H6").Value = 142997.86
$ws.Range("I6").Value = 166814.17
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 500442.51
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = -500329.51
$ws.Range("N6").Value = -526
$ws.Range("H47").Value = 333.75
$ws.Range("I47").Value = 411.66666
$ws.Range("K47").Value = 1234.99998
$ws.Range("M47").Value = -803.9999800000001
$ws.Range("H86").Value = 1214.8572
$ws.Range("I86").Value = 868.3333
$ws.Range("J86").Value = 1474.75
$ws.Range("K86").Value = 2604.9999
$ws.Range("L86").Value = 4424.25
$ws.Range("M86").Value = -1418.9999
$ws.Range("N86").Value = -6796.25
$ws.Range("H89").Value = 1214.8572
$ws.Range("I89").Value = 868.3333
$ws.Range("J89").Value = 1474.75
$ws.Range("K89").Value = 7814.9997
$ws.Range("L89").Value = 13272.75
$ws.Range("M89").Value = -1886.9997
$ws.Range("N89").Value = -25128.75
$ws.Range("H132").Value = 986.5714
$ws.Range("J132").Value = 897.5
$ws.Range("L132").Value = 8077.5
$ws.Range("N132").Value = -13137.5
$ws.Range("H133").Value = 8397.5
$ws.Range("I133").Value = 8397.5
$ws.Range("K133").Value = 25192.5
$ws.Range("M133").Value = -20132.5

$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("M2").Value = 13
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H7").Value = 43039.848
$ws.Range("I7").Value = 54501.8
$ws.Range("K7").Value = 54501.8
$ws.Range("M7").Value = -54389.8
$ws.Range("H8").Value = 43039.848
$ws.Range("I8").Value = 54501.8
$ws.Range("K8").Value = 54501.8
$ws.Range("M8").Value = -54362.8
$ws.Range("H43").Value = 21580.75
$ws.Range("J43").Value = 63125
$ws.Range("L43").Value = 63125
$ws.Range("N43").Value = -63427
$ws.Range("H57").Value = 37999.875
$ws.Range("J57").Value = 44833
$ws.Range("L57").Value = 44833
$ws.Range("N57").Value = -46473
$ws.Range("H70").Value = 3860.5557
$ws.Range("I70").Value = 3998.75
$ws.Range("J70").Value = 3750
$ws.Range("K70").Value = 3998.75
$ws.Range("L70").Value = 3750
$ws.Range("M70").Value = -3728.75
$ws.Range("N70").Value = -4290
$ws.Range("H73").Value = 3860.5557
$ws.Range("I73").Value = 3998.75
$ws.Range("J73").Value = 3750
$ws.Range("K73").Value = 3998.75
$ws.Range("L73").Value = 3750
$ws.Range("M73").Value = -3062.75
$ws.Range("N73").Value = -5622
$ws.Range("H95").Value = 23999
$ws.Range("J95").Value = 23999
$ws.Range("L95").Value = 23999
$ws.Range("N95").Value = -29491
$ws.Range("H107").Value = 3135.625
$ws.Range("J107").Value = 7332.6665
$ws.Range("L107").Value = 7332.6665
$ws.Range("N107").Value = -11172.6665
$ws.Range("H113").Value = 2281.5386
$ws.Range("I113").Value = 2724.111
$ws.Range("J113").Value = 1285.75
$ws.Range("K113").Value = 2724.111
$ws.Range("L113").Value = 1285.75
$ws.Range("M113").Value = -554.1109999999999
$ws.Range("N113").Value = -5625.75
$ws.Range("H122").Value = 3060.8572
$ws.Range("J122").Value = 3203.8
$ws.Range("L122").Value = 9611.400000000001
$ws.Range("N122").Value = -14511.4
$ws.Range("H132").Value = 19610982
$ws.Range("I132").Value = 2779.9167
$ws.Range("K132").Value = 8339.750100000001
$ws.Range("M132").Value = -5809.750100000001

$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 4747.5
$ws.Range("I7").Value = 4000
$ws.Range("J7").Value = 4996.6665
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 4996.6665
$ws.Range("M7").Value = -3888
$ws.Range("N7").Value = -5220.6665
$ws.Range("H40").Value = 35718350
$ws.Range("I40").Value = 52634764
$ws.Range("K40").Value = 52634764
$ws.Range("M40").Value = -52634628
$ws.Range("H42").Value = 39991
$ws.Range("J42").Value = 39991
$ws.Range("L42").Value = 39991
$ws.Range("N42").Value = -41117
$ws.Range("H49").Value = 39991
$ws.Range("J49").Value = 39991
$ws.Range("L49").Value = 39991
$ws.Range("N49").Value = -40285
$ws.Range("H55").Value = 850.0417
$ws.Range("I55").Value = 407.5
$ws.Range("J55").Value = 1469.6
$ws.Range("K55").Value = 407.5
$ws.Range("L55").Value = 1469.6
$ws.Range("M55").Value = -234.5
$ws.Range("N55").Value = -1815.6
$ws.Range("H100").Value = 2034.3334
$ws.Range("I100").Value = 1001.5
$ws.Range("K100").Value = 1001.5
$ws.Range("M100").Value = -460.5
$ws.Range("H126").Value = 4747.5
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 4996.6665
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 14989.9995
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -19929.9995

$ws = $wb.Worksheets.Item(8)
$ws.Range("H40").Value = 24900
$ws.Range("J40").Value = 24900
$ws.Range("L40").Value = 24900
$ws.Range("N40").Value = -25198
$ws.Range("H54").Value = 47806
$ws.Range("J54").Value = 47806
$ws.Range("L54").Value = 47806
$ws.Range("N54").Value = -48846
$ws.Range("H74").Value = 22164
$ws.Range("I74").Value = 20662
$ws.Range("K74").Value = 20662
$ws.Range("M74").Value = -19726
$ws.Range("H77").Value = 22164
$ws.Range("I77").Value = 20662
$ws.Range("K77").Value = 61986
$ws.Range("M77").Value = -57306
$ws.Range("H122").Value = 1157.7916
$ws.Range("J122").Value = 467.33334
$ws.Range("L122").Value = 1402.00002
$ws.Range("N122").Value = -6302.000019999999
